$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose permission-check result moved from "Vào được" (granted) to
# "Đã hủy" (revoked) for columns B:E (Khách / Nhân viên / Chi nhánh / Admin).
$rows = @(68, 69, 70, 71, 72, 74)

# Use the format already sitting on D75 (black Tahoma text, no fill) as the
# donor for the font, then paint the yellow highlight back on top - this is
# exactly the look ("Đã hủy" rows) used elsewhere in the sheet.
$ws.Range("D75").Copy() | Out-Null

foreach ($r in $rows) {
    $rng = $ws.Range("B" + $r + ":E" + $r)
    $rng.PasteSpecial(-4122) | Out-Null
    $rng.Interior.Color = 65535
    $rng.Value = "Đã hủy"
}

$excel.CutCopyMode = 0

# Restore selection to the row that was just updated.
$ws.Range("B74:E74").Select() | Out-Null
